# Add a new "Requirement" worksheet after the existing sheets (DMCAR, Mapping)
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$reqSheet = $wb.Worksheets.Add($null, $lastSheet)
$reqSheet.Name = "Requirement"

# Populate the header row for the new Requirement sheet
$reqSheet.Cells.Item(1, 1).Value = "Project"
$reqSheet.Cells.Item(1, 2).Value = "PrimaryStakeholder"
$reqSheet.Cells.Item(1, 3).Value = "PrimaryStakeholderDescription"
$reqSheet.Cells.Item(1, 4).Value = "RequirementId"
$reqSheet.Cells.Item(1, 5).Value = "Name"
$reqSheet.Cells.Item(1, 6).Value = "Description"

# Give the new columns a sensible custom width (close to the authored layout)
$reqSheet.Columns.Item(1).ColumnWidth = 15.8203125
$reqSheet.Columns.Item(2).ColumnWidth = 19.7265625
$reqSheet.Columns.Item(3).ColumnWidth = 28.7109375
$reqSheet.Columns.Item(4).ColumnWidth = 17.578125
$reqSheet.Columns.Item(5).ColumnWidth = 13.0859375
$reqSheet.Columns.Item(6).ColumnWidth = 15.0390625
$reqSheet.Columns.Item(7).ColumnWidth = 15.4296875

# Restore/update the view state (selection) on the pre-existing sheets
$dmcar = $wb.Worksheets.Item("DMCAR")
[void]$dmcar.Activate()
[void]$dmcar.Range("B8").Select()
$excel.ActiveWindow.DisplayGridlines = $true

$mapping = $wb.Worksheets.Item("Mapping")
[void]$mapping.Activate()
[void]$mapping.Range("H4").Select()
$excel.ActiveWindow.DisplayGridlines = $true

# Make the new Requirement sheet the active/selected tab (activeTab + tabSelected)
[void]$reqSheet.Activate()
[void]$reqSheet.Range("B8").Select()
$excel.ActiveWindow.DisplayGridlines = $true
